# Workbook starts with a single sheet ("Sheet1") that already holds the
# "mean" statistics. The edit adds five more sheets (one per TSP instance)
# in front of it, each with the same layout/style but its own numbers, and
# renames the original sheet to "mean" so it becomes the last tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the existing sheet 5 times, always appending after the most
# recently created copy so the tab order stays left-to-right.
for ($i = 1; $i -le 5; $i++) {
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $last.Copy($null, $last)
}

# Rename the six tabs into their final order.
$names = @("eil51", "berlin52", "pr136", "pr226", "d198", "mean")
for ($i = 1; $i -le 6; $i++) {
    $wb.Worksheets.Item($i).Name = $names[$i - 1]
}

# Per-sheet statistics: row 2 = "auc", row 3 = "min", columns B:E = random,
# bayesian, forest, gradient. "mean" keeps the values that were already in
# the workbook, so it isn't touched here.
$data = @{
    "eil51"    = @{ row2 = @(0.8304121401559936, 0.899677890294611, 0.8188253837210422, 0.601256411976037);
                    row3 = @(0.02623362030768031, 0.03391590397004692, 0.01968897357784272, 0.02395574176995417) }
    "berlin52" = @{ row2 = @(0.3470775386038567, 0.6504245201931586, 0.2266698501351931, 0.1590209092862858);
                    row3 = @(0.0007518913878295021, 0.01900986082063078, 0.004513743086632056, 0.003085628744123837) }
    "pr136"    = @{ row2 = @(2.265570492570207, 1.849572900070808, 1.809140861395097, 1.947604744127015);
                    row3 = @(0.09676929003141378, 0.09267605715685522, 0.07540059859748173, 0.07542778875383883) }
    "pr226"    = @{ row2 = @(1.837164294898639, 1.394382597025814, 1.54691092308816, 1.497230361062779);
                    row3 = @(0.07454252915968133, 0.07226467000338908, 0.07155621950623137, 0.04211928288038776) }
    "d198"     = @{ row2 = @(2.138673722391764, 1.697629530878448, 1.940128028161632, 1.745387991281444);
                    row3 = @(0.09220523859069986, 0.08204043426292522, 0.08460608390260768, 0.08085682164963927) }
}

$cols = @("B", "C", "D", "E")

foreach ($sheetName in $data.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $vals2 = $data[$sheetName].row2
    $vals3 = $data[$sheetName].row3
    for ($c = 0; $c -lt 4; $c++) {
        $ws.Range("$($cols[$c])2").Value = $vals2[$c]
        $ws.Range("$($cols[$c])3").Value = $vals3[$c]
    }
}

$wb.Worksheets.Item("eil51").Activate()
